$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9958560466766357
$ws.Range("B1").Value = 2.10805869102478
$ws.Range("C1").Value = 6.946327686309814
$ws.Range("D1").Value = 2.146793603897095
$ws.Range("E1").Value = 1.37790310382843
